$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the series. It belongs right
# after the header/date-18 row, so insert a fresh row at position 19 and
# push all the existing data rows (old 19..63) down to (20..64).
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new observation.
$ws.Cells.Item(19, 1).Value = 9
$ws.Cells.Item(19, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(19, 3).Value = "Metropolitana"
$ws.Cells.Item(19, 4).Value = 44804
$ws.Cells.Item(19, 5).Value = 13
$ws.Cells.Item(19, 6).Value = 100112035
$ws.Cells.Item(19, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 35
$ws.Cells.Item(19, 11).Value = 19000
$ws.Cells.Item(19, 12).Value = 20000
$ws.Cells.Item(19, 13).Value = 19714
$ws.Cells.Item(19, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(19, 15).Value = "Hijuelas"
$ws.Cells.Item(19, 16).Value = 1314
$ws.Cells.Item(19, 17).Value = 15
$ws.Cells.Item(19, 18).Value = "Hortaliza"
